$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the "contingency" cell at A2: value 0, bold font, thin box border,
# centered horizontally / top vertically aligned.
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Borders.Weight = 2
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

# B1 gets the same formatting (value 0) - copy the formatting from A2 so both
# cells share a single style entry instead of each accumulating their own.
$ws.Range("B1").Value = 0
$ws.Range("A2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B2 holds the label, stored as a shared string, with default formatting.
$ws.Range("B2").Value = "disconnected_elements"
